# Scheduled market-data refresh: update computed Leve profit columns (H:N)
# on each profession sheet to reflect newly-fetched Universalis prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 334.9  # H4
$ws.Cells.Item(4, 9).Value = 257.125  # I4
$ws.Cells.Item(4, 10).Value = 646  # J4
$ws.Cells.Item(4, 11).Value = 257.125  # K4
$ws.Cells.Item(4, 12).Value = 646  # L4
$ws.Cells.Item(4, 13).Value = -143.125  # M4
$ws.Cells.Item(4, 14).Value = -874  # N4
$ws.Cells.Item(28, 8).Value = 824.8570999999999  # H28
$ws.Cells.Item(28, 9).Value = 432.3125  # I28
$ws.Cells.Item(28, 10).Value = 2081  # J28
$ws.Cells.Item(28, 11).Value = 432.3125  # K28
$ws.Cells.Item(28, 12).Value = 2081  # L28
$ws.Cells.Item(28, 13).Value = 52.6875  # M28
$ws.Cells.Item(28, 14).Value = -3051  # N28
$ws.Cells.Item(86, 8).Value = 4750.375  # H86
$ws.Cells.Item(86, 9).Value = 4501  # I86
$ws.Cells.Item(86, 10).Value = 4900  # J86
$ws.Cells.Item(86, 11).Value = 4501  # K86
$ws.Cells.Item(86, 12).Value = 4900  # L86
$ws.Cells.Item(86, 13).Value = -3378  # M86
$ws.Cells.Item(86, 14).Value = -7146  # N86
$ws.Cells.Item(89, 8).Value = 4750.375  # H89
$ws.Cells.Item(89, 9).Value = 4501  # I89
$ws.Cells.Item(89, 10).Value = 4900  # J89
$ws.Cells.Item(89, 11).Value = 22505  # K89
$ws.Cells.Item(89, 12).Value = 24500  # L89
$ws.Cells.Item(89, 13).Value = -16889  # M89
$ws.Cells.Item(89, 14).Value = -35732  # N89
$ws.Cells.Item(132, 8).Value = 4388432  # H132
$ws.Cells.Item(132, 9).Value = 2325.653  # I132
$ws.Cells.Item(132, 10).Value = 31253334  # J132
$ws.Cells.Item(132, 11).Value = 6976.958999999999  # K132
$ws.Cells.Item(132, 12).Value = 93760002  # L132
$ws.Cells.Item(132, 13).Value = -4446.958999999999  # M132
$ws.Cells.Item(132, 14).Value = -93765062  # N132
$ws.Cells.Item(141, 8).Value = 1328.56  # H141
$ws.Cells.Item(141, 9).Value = 1183.2174  # I141
$ws.Cells.Item(141, 11).Value = 3549.6522  # K141
$ws.Cells.Item(141, 13).Value = 1630.3478  # M141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6265.88  # H32
$ws.Cells.Item(32, 9).Value = 5280.1377  # I32
$ws.Cells.Item(32, 10).Value = 12862.77  # J32
$ws.Cells.Item(32, 11).Value = 5280.1377  # K32
$ws.Cells.Item(32, 12).Value = 12862.77  # L32
$ws.Cells.Item(32, 13).Value = -4993.1377  # M32
$ws.Cells.Item(32, 14).Value = -13436.77  # N32
$ws.Cells.Item(61, 8).Value = 2140.75  # H61
$ws.Cells.Item(61, 9).Value = 2189.4473  # I61
$ws.Cells.Item(61, 11).Value = 2189.4473  # K61
$ws.Cells.Item(61, 13).Value = -1977.4473  # M61
$ws.Cells.Item(74, 8).Value = 1461.4857  # H74
$ws.Cells.Item(74, 9).Value = 1160.7084  # I74
$ws.Cells.Item(74, 11).Value = 1160.7084  # K74
$ws.Cells.Item(74, 13).Value = -286.7084  # M74
$ws.Cells.Item(77, 8).Value = 1461.4857  # H77
$ws.Cells.Item(77, 9).Value = 1160.7084  # I77
$ws.Cells.Item(77, 11).Value = 5803.541999999999  # K77
$ws.Cells.Item(77, 13).Value = -1435.541999999999  # M77
$ws.Cells.Item(102, 8).Value = 1500  # H102
$ws.Cells.Item(102, 9).Value = 0  # I102
$ws.Cells.Item(102, 10).Value = 1500  # J102
$ws.Cells.Item(102, 11).Value = 0  # K102
$ws.Cells.Item(102, 12).Value = 1500  # L102
$ws.Cells.Item(102, 13).ClearContents()  # M102
$ws.Cells.Item(102, 14).Value = -4744  # N102
$ws.Cells.Item(136, 8).Value = 2140.75  # H136
$ws.Cells.Item(136, 9).Value = 2189.4473  # I136
$ws.Cells.Item(136, 11).Value = 6568.341899999999  # K136
$ws.Cells.Item(136, 13).Value = -4018.341899999999  # M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(75, 8).Value = 16673.666  # H75
$ws.Cells.Item(75, 9).Value = 3225.8333  # I75
$ws.Cells.Item(75, 10).Value = 43569.332  # J75
$ws.Cells.Item(75, 11).Value = 3225.8333  # K75
$ws.Cells.Item(75, 12).Value = 43569.332  # L75
$ws.Cells.Item(75, 13).Value = -2289.8333  # M75
$ws.Cells.Item(75, 14).Value = -45441.332  # N75
$ws.Cells.Item(78, 8).Value = 16673.666  # H78
$ws.Cells.Item(78, 9).Value = 3225.8333  # I78
$ws.Cells.Item(78, 10).Value = 43569.332  # J78
$ws.Cells.Item(78, 11).Value = 9677.499899999999  # K78
$ws.Cells.Item(78, 12).Value = 130707.996  # L78
$ws.Cells.Item(78, 13).Value = -4997.499899999999  # M78
$ws.Cells.Item(78, 14).Value = -140067.996  # N78
$ws.Cells.Item(103, 8).Value = 40528.5  # H103
$ws.Cells.Item(103, 10).Value = 40528.5  # J103
$ws.Cells.Item(103, 12).Value = 40528.5  # L103
$ws.Cells.Item(103, 14).Value = -42872.5  # N103

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2167.1292  # H132
$ws.Cells.Item(132, 9).Value = 1712.5  # I132
$ws.Cells.Item(132, 11).Value = 5137.5  # K132
$ws.Cells.Item(132, 13).Value = -2607.5  # M132
$ws.Cells.Item(134, 8).Value = 701763.7  # H134
$ws.Cells.Item(134, 9).Value = 1651.4138  # I134
$ws.Cells.Item(134, 10).Value = 4762415  # J134
$ws.Cells.Item(134, 11).Value = 4954.2414  # K134
$ws.Cells.Item(134, 12).Value = 14287245  # L134
$ws.Cells.Item(134, 13).Value = -2419.2414  # M134
$ws.Cells.Item(134, 14).Value = -14292315  # N134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1711.5714  # H5
$ws.Cells.Item(5, 9).Value = 314.81818  # I5
$ws.Cells.Item(5, 11).Value = 944.45454  # K5
$ws.Cells.Item(5, 13).Value = -832.45454  # M5
$ws.Cells.Item(92, 8).Value = 1239.6364  # H92
$ws.Cells.Item(92, 9).Value = 300  # I92
$ws.Cells.Item(92, 10).Value = 1333.6  # J92
$ws.Cells.Item(92, 11).Value = 900  # K92
$ws.Cells.Item(92, 12).Value = 4000.8  # L92
$ws.Cells.Item(92, 13).Value = 348  # M92
$ws.Cells.Item(92, 14).Value = -6496.799999999999  # N92
$ws.Cells.Item(123, 8).Value = 5940.864  # H123
$ws.Cells.Item(123, 9).Value = 2333.3333  # I123
$ws.Cells.Item(123, 11).Value = 6999.999899999999  # K123
$ws.Cells.Item(123, 13).Value = -4549.999899999999  # M123
$ws.Cells.Item(135, 8).Value = 1711.5714  # H135
$ws.Cells.Item(135, 9).Value = 314.81818  # I135
$ws.Cells.Item(135, 11).Value = 2833.36362  # K135
$ws.Cells.Item(135, 13).Value = -298.3636200000001  # M135

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 82  # H2
$ws.Cells.Item(2, 9).Value = 46.75  # I2
$ws.Cells.Item(2, 10).Value = 138.4  # J2
$ws.Cells.Item(2, 11).Value = 46.75  # K2
$ws.Cells.Item(2, 12).Value = 138.4  # L2
$ws.Cells.Item(2, 13).Value = 66.25  # M2
$ws.Cells.Item(2, 14).Value = -364.4  # N2
$ws.Cells.Item(64, 8).Value = 15000  # H64
$ws.Cells.Item(64, 10).Value = 15000  # J64
$ws.Cells.Item(64, 12).Value = 15000  # L64
$ws.Cells.Item(64, 14).Value = -15496  # N64
$ws.Cells.Item(67, 8).Value = 15000  # H67
$ws.Cells.Item(67, 10).Value = 15000  # J67
$ws.Cells.Item(67, 12).Value = 15000  # L67
$ws.Cells.Item(67, 14).Value = -16716  # N67
$ws.Cells.Item(80, 8).Value = 15875695  # H80
$ws.Cells.Item(80, 9).Value = 20835706  # I80
$ws.Cells.Item(80, 10).Value = 3661.2  # J80
$ws.Cells.Item(80, 11).Value = 20835706  # K80
$ws.Cells.Item(80, 12).Value = 3661.2  # L80
$ws.Cells.Item(80, 13).Value = -20834708  # M80
$ws.Cells.Item(80, 14).Value = -5657.2  # N80
$ws.Cells.Item(83, 8).Value = 15875695  # H83
$ws.Cells.Item(83, 9).Value = 20835706  # I83
$ws.Cells.Item(83, 10).Value = 3661.2  # J83
$ws.Cells.Item(83, 11).Value = 104178530  # K83
$ws.Cells.Item(83, 12).Value = 18306  # L83
$ws.Cells.Item(83, 13).Value = -104173538  # M83
$ws.Cells.Item(83, 14).Value = -28290  # N83
$ws.Cells.Item(102, 8).Value = 2588.5151  # H102
$ws.Cells.Item(102, 9).Value = 3059  # I102
$ws.Cells.Item(102, 10).Value = 1647.5454  # J102
$ws.Cells.Item(102, 11).Value = 3059  # K102
$ws.Cells.Item(102, 12).Value = 1647.5454  # L102
$ws.Cells.Item(102, 13).Value = -1437  # M102
$ws.Cells.Item(102, 14).Value = -4891.5454  # N102
$ws.Cells.Item(113, 8).Value = 92206.09  # H113
$ws.Cells.Item(113, 9).Value = 101176.7  # I113
$ws.Cells.Item(113, 10).Value = 2500  # J113
$ws.Cells.Item(113, 11).Value = 101176.7  # K113
$ws.Cells.Item(113, 12).Value = 2500  # L113
$ws.Cells.Item(113, 13).Value = -99006.7  # M113
$ws.Cells.Item(113, 14).Value = -6840  # N113
$ws.Cells.Item(132, 8).Value = 5136.8647  # H132
$ws.Cells.Item(132, 9).Value = 5944.4443  # I132
$ws.Cells.Item(132, 10).Value = 2956.4  # J132
$ws.Cells.Item(132, 11).Value = 17833.3329  # K132
$ws.Cells.Item(132, 12).Value = 8869.200000000001  # L132
$ws.Cells.Item(132, 13).Value = -15303.3329  # M132
$ws.Cells.Item(132, 14).Value = -13929.2  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1084.1578  # H22
$ws.Cells.Item(22, 9).Value = 561.375  # I22
$ws.Cells.Item(22, 10).Value = 1464.3636  # J22
$ws.Cells.Item(22, 11).Value = 561.375  # K22
$ws.Cells.Item(22, 12).Value = 1464.3636  # L22
$ws.Cells.Item(22, 13).Value = -266.375  # M22
$ws.Cells.Item(22, 14).Value = -2054.3636  # N22
$ws.Cells.Item(27, 8).Value = 1084.1578  # H27
$ws.Cells.Item(27, 9).Value = 561.375  # I27
$ws.Cells.Item(27, 10).Value = 1464.3636  # J27
$ws.Cells.Item(27, 11).Value = 561.375  # K27
$ws.Cells.Item(27, 12).Value = 1464.3636  # L27
$ws.Cells.Item(27, 13).Value = -454.375  # M27
$ws.Cells.Item(27, 14).Value = -1678.3636  # N27
$ws.Cells.Item(132, 8).Value = 7944644  # H132
$ws.Cells.Item(132, 9).Value = 6355.375  # I132
$ws.Cells.Item(132, 10).Value = 33347166  # J132
$ws.Cells.Item(132, 11).Value = 19066.125  # K132
$ws.Cells.Item(132, 12).Value = 100041498  # L132
$ws.Cells.Item(132, 13).Value = -16536.125  # M132
$ws.Cells.Item(132, 14).Value = -100046558  # N132
$ws.Cells.Item(136, 8).Value = 5215.4326  # H136
$ws.Cells.Item(136, 9).Value = 2271.394  # I136
$ws.Cells.Item(136, 11).Value = 6814.181999999999  # K136
$ws.Cells.Item(136, 13).Value = -4264.181999999999  # M136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 944.44446  # H81
$ws.Cells.Item(81, 10).Value = 1600  # J81
$ws.Cells.Item(81, 12).Value = 3200  # L81
$ws.Cells.Item(81, 14).Value = -5322  # N81
$ws.Cells.Item(84, 8).Value = 944.44446  # H84
$ws.Cells.Item(84, 10).Value = 1600  # J84
$ws.Cells.Item(84, 12).Value = 16000  # L84
$ws.Cells.Item(84, 14).Value = -26608  # N84
$ws.Cells.Item(107, 8).Value = 1450.8889  # H107
$ws.Cells.Item(107, 9).Value = 1519.75  # I107
$ws.Cells.Item(107, 10).Value = 900  # J107
$ws.Cells.Item(107, 11).Value = 4559.25  # K107
$ws.Cells.Item(107, 12).Value = 2700  # L107
$ws.Cells.Item(107, 13).Value = -2639.25  # M107
$ws.Cells.Item(107, 14).Value = -6540  # N107
$ws.Cells.Item(113, 8).Value = 797.4666999999999  # H113
$ws.Cells.Item(113, 9).Value = 1008.5455  # I113
$ws.Cells.Item(113, 10).Value = 217  # J113
$ws.Cells.Item(113, 11).Value = 3025.6365  # K113
$ws.Cells.Item(113, 12).Value = 651  # L113
$ws.Cells.Item(113, 13).Value = -855.6364999999996  # M113
$ws.Cells.Item(113, 14).Value = -4991  # N113
$ws.Cells.Item(122, 8).Value = 2571.0881  # H122
$ws.Cells.Item(122, 9).Value = 2480.4666  # I122
$ws.Cells.Item(122, 10).Value = 3250.75  # J122
$ws.Cells.Item(122, 11).Value = 7441.399800000001  # K122
$ws.Cells.Item(122, 12).Value = 9752.25  # L122
$ws.Cells.Item(122, 13).Value = -4991.399800000001  # M122
$ws.Cells.Item(122, 14).Value = -14652.25  # N122
